$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g18.3")

$ws.Range("B22").Value = 22.63893950014812
$ws.Range("C22").Value = 11.48384476288428
$ws.Range("D22").Value = 13.51292263323663

$ws.Range("B23").Value = 23.09934200987477
$ws.Range("C23").Value = 11.82184518019592
$ws.Range("D23").Value = 13.74167909528762

$ws.Range("B24").Value = 23.5237729586581
$ws.Range("C24").Value = 12.14811097263668
$ws.Range("D24").Value = 13.98911715807688

$ws.Range("B25").Value = 23.70743431183153
$ws.Range("C25").Value = 12.43147182202182
$ws.Range("D25").Value = 14.01091629647116

$ws.Range("B26").Value = 23.88572363975156
$ws.Range("C26").Value = 12.58166955787716
$ws.Range("D26").Value = 14.1877331400001
